$d = $word.ActiveDocument

# Change 1: "mes(es)" -> "año(s)" inside the intro paragraph
$d.Content.Find.Execute("para el(los) mes(es) {{AÑO}}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "para el(los) año(s) {{AÑO}}", 2)

# Change 2: update the cached date field text
$d.Content.Find.Execute("1 de noviembre de 2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "24 de noviembre de 2024", 2)
